$wb = $excel.ActiveWorkbook

# --- Update the formulas / cached values on the "air-dry-2019-whc" sheet ---
# Column A rows 2-13 compare WHC!A(row+3) to "Field" and return either
# "control" or "air-dry" (previously "control (2019)" / "air-dry (2019)").
$ws = $wb.Worksheets.Item("air-dry-2019-whc")

for ($r = 2; $r -le 13; $r++) {
    $whcRow = $r + 3
    $ws.Range("A$r").Formula = "=IF(WHC!A$whcRow=""Field"",""control"",""air-dry"")"
}

# --- Update sheet view / selection state ---

# "meta" sheet should no longer be the selected tab.
$metaWs = $wb.Worksheets.Item("meta")
$metaWs.Select()

# "air-dry-2019-whc" becomes the active/selected sheet, with a new
# selection of F6 (instead of A8:A13).
$ws.Activate()
$ws.Range("F6").Select()
